$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1067.6957
$ws.Range("I107").Value = 897.3333
$ws.Range("K107").Value = 897.3333
$ws.Range("M107").Value = 1022.6667

$ws.Range("H125").Value = 407.3
$ws.Range("I125").Value = 366.4
$ws.Range("J125").Value = 448.2
$ws.Range("K125").Value = 3297.6
$ws.Range("L125").Value = 4033.8
$ws.Range("M125").Value = -837.5999999999999
$ws.Range("N125").Value = -8953.799999999999

$ws.Range("H129").Value = 1541.8286
$ws.Range("J129").Value = 1667.7188
$ws.Range("L129").Value = 5003.1564
$ws.Range("N129").Value = -15003.1564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1998
$ws.Range("I19").Value = 1998
$ws.Range("K19").Value = 1998
$ws.Range("M19").Value = -1769

$ws.Range("H45").Value = 2472.2
$ws.Range("I45").Value = 2159.96
$ws.Range("K45").Value = 2159.96
$ws.Range("M45").Value = -1782.96

$ws.Range("H97").Value = 995
$ws.Range("I97").Value = 995
$ws.Range("K97").Value = 995
$ws.Range("M97").Value = -499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2685.2942
$ws.Range("I20").Value = 2242.3076
$ws.Range("J20").Value = 4125
$ws.Range("K20").Value = 2242.3076
$ws.Range("L20").Value = 4125
$ws.Range("M20").Value = -1995.3076
$ws.Range("N20").Value = -4619

$ws.Range("H94").Value = 697.0769
$ws.Range("I94").Value = 562.9286
$ws.Range("J94").Value = 853.5833
$ws.Range("K94").Value = 562.9286
$ws.Range("L94").Value = 853.5833
$ws.Range("M94").Value = -111.9286
$ws.Range("N94").Value = -1755.5833

$ws.Range("H105").Value = 1430637.9
$ws.Range("I105").Value = 1386.5834
$ws.Range("K105").Value = 1386.5834
$ws.Range("M105").Value = 360.4166

$ws.Range("H107").Value = 1687.8649
$ws.Range("I107").Value = 1380.7
$ws.Range("K107").Value = 1380.7
$ws.Range("M107").Value = 539.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1135.7778
$ws.Range("I16").Value = 1237
$ws.Range("K16").Value = 1237
$ws.Range("M16").Value = -950

$ws.Range("H94").Value = 6304
$ws.Range("I94").Value = 4081
$ws.Range("J94").Value = 10750
$ws.Range("K94").Value = 4081
$ws.Range("L94").Value = 10750
$ws.Range("M94").Value = -3630
$ws.Range("N94").Value = -11652

$ws.Range("H105").Value = 4700
$ws.Range("I105").Value = 5550
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 5550
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -3803
$ws.Range("N105").Value = -6494

$ws.Range("H107").Value = 1373.15
$ws.Range("I107").Value = 505.85715
$ws.Range("J107").Value = 1840.1538
$ws.Range("K107").Value = 505.85715
$ws.Range("L107").Value = 1840.1538
$ws.Range("M107").Value = 1414.14285
$ws.Range("N107").Value = -5680.1538

$ws.Range("H113").Value = 1135.7778
$ws.Range("I113").Value = 1237
$ws.Range("K113").Value = 1237
$ws.Range("M113").Value = 933

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 729.9400000000001
$ws.Range("I131").Value = 380.55554
$ws.Range("J131").Value = 764.4945
$ws.Range("K131").Value = 1141.66662
$ws.Range("L131").Value = 2293.4835
$ws.Range("M131").Value = 3898.33338
$ws.Range("N131").Value = -12373.4835

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 159461.53
$ws.Range("I43").Value = 6083.3335
$ws.Range("K43").Value = 6083.3335
$ws.Range("M43").Value = -5932.3335

$ws.Range("H46").Value = 34800
$ws.Range("J46").Value = 34800
$ws.Range("L46").Value = 34800
$ws.Range("N46").Value = -35112

$ws.Range("H57").Value = 29227.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 29227.5
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 29227.5
$ws.Range("N57").Value = -30867.5
$ws.Range("M57").ClearContents()

$ws.Range("H70").Value = 3129387.8
$ws.Range("I70").Value = 4529.7
$ws.Range("J70").Value = 6254246
$ws.Range("K70").Value = 4529.7
$ws.Range("L70").Value = 6254246
$ws.Range("M70").Value = -4259.7
$ws.Range("N70").Value = -6254786

$ws.Range("H73").Value = 3129387.8
$ws.Range("I73").Value = 4529.7
$ws.Range("J73").Value = 6254246
$ws.Range("K73").Value = 4529.7
$ws.Range("L73").Value = 6254246
$ws.Range("M73").Value = -3593.7
$ws.Range("N73").Value = -6256118

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4400
$ws.Range("I7").Value = 3575
$ws.Range("J7").Value = 7700
$ws.Range("K7").Value = 3575
$ws.Range("L7").Value = 7700
$ws.Range("M7").Value = -3463
$ws.Range("N7").Value = -7924

$ws.Range("H40").Value = 4418
$ws.Range("I40").Value = 4418
$ws.Range("K40").Value = 4418
$ws.Range("M40").Value = -4282

$ws.Range("H64").Value = 38000
$ws.Range("J64").Value = 38000
$ws.Range("L64").Value = 38000
$ws.Range("N64").Value = -38450

$ws.Range("H67").Value = 38000
$ws.Range("J67").Value = 38000
$ws.Range("L67").Value = 38000
$ws.Range("N67").Value = -39560

$ws.Range("H93").Value = 2797.4375
$ws.Range("I93").Value = 2746.5833
$ws.Range("J93").Value = 2950
$ws.Range("K93").Value = 2746.5833
$ws.Range("L93").Value = 2950
$ws.Range("M93").Value = -1498.5833
$ws.Range("N93").Value = -5446

$ws.Range("H100").Value = 2306.375
$ws.Range("I100").Value = 1701
$ws.Range("J100").Value = 2581.5454
$ws.Range("K100").Value = 1701
$ws.Range("L100").Value = 2581.5454
$ws.Range("M100").Value = -1160
$ws.Range("N100").Value = -3663.5454

$ws.Range("H126").Value = 4400
$ws.Range("I126").Value = 3575
$ws.Range("J126").Value = 7700
$ws.Range("K126").Value = 10725
$ws.Range("L126").Value = 23100
$ws.Range("M126").Value = -8255
$ws.Range("N126").Value = -28040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 17361.25
$ws.Range("J70").Value = 19783.334
$ws.Range("L70").Value = 19783.334
$ws.Range("N70").Value = -20413.334

$ws.Range("H73").Value = 17361.25
$ws.Range("J73").Value = 19783.334
$ws.Range("L73").Value = 19783.334
$ws.Range("N73").Value = -21967.334

$ws.Range("H126").Value = 2496.5881
$ws.Range("I126").Value = 2125
$ws.Range("K126").Value = 6375
$ws.Range("M126").Value = -3905

$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 20000
$ws.Range("N131").Value = -30080
